# ------------------------------------------------------------------
# Refresh the crypto ranking table (Coin / Link / Price / Volume(1h)).
#
# 1) Rows 2-33 keep their ranking position; only Price (D) and
#    Volume(1h) (E) are refreshed with the latest quote.
# 2) A new coin ('Frax') now ranks in at row 34, so Coin/Link/Price/
#    Volume (columns B:E) for every following coin shift down by one
#    row. Column A (the 0-based rank index) is left untouched, since
#    it already lists 0..49 top-to-bottom regardless of which coin
#    occupies the row.
#
# Price values are entered with a leading apostrophe so Excel stores
# them as literal text (matching the source feed's formatting, e.g.
# '26.726.07' or '0.9980') instead of silently reinterpreting them
# as numbers and dropping trailing zeros / changing notation.
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: refresh Price / Volume(1h) for rows 2-33 ---
$ws.Cells.Item(2, 4).Value = '''26.726.07'
$ws.Cells.Item(2, 5).Value = '  +1.33%  '
$ws.Cells.Item(3, 4).Value = '''1.733.86'
$ws.Cells.Item(3, 5).Value = '  +0.72%  '
$ws.Cells.Item(4, 4).Value = '''0.9977'
$ws.Cells.Item(4, 5).Value = '  -0.27%  '
$ws.Cells.Item(5, 4).Value = '''242.52'
$ws.Cells.Item(5, 5).Value = '  -0.64%  '
$ws.Cells.Item(6, 4).Value = '''0.9980'
$ws.Cells.Item(6, 5).Value = '  -0.30%  '
$ws.Cells.Item(7, 5).Value = '  +0.91%  '
$ws.Cells.Item(8, 5).Value = '  +0.56%  '
$ws.Cells.Item(9, 4).Value = '''0.06221'
$ws.Cells.Item(9, 5).Value = '  +0.27%  '
$ws.Cells.Item(10, 4).Value = '''1.727.01'
$ws.Cells.Item(10, 5).Value = '  +0.33%  '
$ws.Cells.Item(11, 4).Value = '''15.94'
$ws.Cells.Item(11, 5).Value = '  +3.44%  '
$ws.Cells.Item(12, 4).Value = '''0.06977'
$ws.Cells.Item(12, 5).Value = '  -0.52%  '
$ws.Cells.Item(13, 4).Value = '''0.6113'
$ws.Cells.Item(13, 5).Value = '  +2.41%  '
$ws.Cells.Item(14, 4).Value = '''4.506'
$ws.Cells.Item(14, 5).Value = '  -0.55%  '
$ws.Cells.Item(15, 4).Value = '''77.24'
$ws.Cells.Item(15, 5).Value = '  +0.15%  '
$ws.Cells.Item(16, 4).Value = '''0.9987'
$ws.Cells.Item(16, 5).Value = '  -0.23%  '
$ws.Cells.Item(17, 4).Value = '''26.519.69'
$ws.Cells.Item(17, 5).Value = '  +0.50%  '
$ws.Cells.Item(18, 4).Value = '''0.9981'
$ws.Cells.Item(18, 5).Value = '  -0.22%  '
$ws.Cells.Item(19, 4).Value = '''0.000007204'
$ws.Cells.Item(19, 5).Value = '  +0.05%  '
$ws.Cells.Item(20, 4).Value = '''11.42'
$ws.Cells.Item(20, 5).Value = '  +0.65%  '
$ws.Cells.Item(21, 4).Value = '''1.948.92'
$ws.Cells.Item(21, 5).Value = '  -0.11%  '
$ws.Cells.Item(22, 4).Value = '''4.474'
$ws.Cells.Item(22, 5).Value = '  -0.01%  '
$ws.Cells.Item(23, 4).Value = '''8.570'
$ws.Cells.Item(23, 5).Value = '  +0.28%  '
$ws.Cells.Item(24, 4).Value = '''5.101'
$ws.Cells.Item(24, 5).Value = '  -1.27%  '
$ws.Cells.Item(25, 4).Value = '''138.44'
$ws.Cells.Item(25, 5).Value = '  +0.79%  '
$ws.Cells.Item(26, 4).Value = '''15.36'
$ws.Cells.Item(27, 4).Value = '''1.771'
$ws.Cells.Item(27, 5).Value = '  +3.51%  '
$ws.Cells.Item(28, 4).Value = '''1.385'
$ws.Cells.Item(28, 5).Value = '  -1.79%  '
$ws.Cells.Item(29, 5).Value = '  -0.56%  '
$ws.Cells.Item(30, 4).Value = '''3.937'
$ws.Cells.Item(30, 5).Value = '  -0.43%  '
$ws.Cells.Item(31, 4).Value = '''0.07985'
$ws.Cells.Item(32, 4).Value = '''3.669'
$ws.Cells.Item(32, 5).Value = '  -0.04%  '
$ws.Cells.Item(33, 4).Value = '''0.04487'
$ws.Cells.Item(33, 5).Value = '  -0.68%  '

# --- Step 2: push rows 34-50 (Coin/Link/Price/Volume only) down to 35-51 ---
# Walk bottom-up so each row is read before it gets overwritten; this
# copies cell contents verbatim (no retyping), so existing text stays text.
for ($r = 50; $r -ge 34; $r--) {
    $dest = $r + 1
    $ws.Range("B$r`:E$r").Copy($ws.Range("B$dest`:E$dest"))
}

# --- Step 3: write the newly-ranked coin into row 34 ---
$ws.Cells.Item(34, 2).Value = 'Frax'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(34, 4).Value = '''0.9976'
$ws.Cells.Item(34, 5).Value = '  -0.25%  '

# --- Step 4: refresh Price / Volume(1h) for the rows that shifted down (35-51) ---
$ws.Cells.Item(35, 4).Value = '''2.608'
$ws.Cells.Item(35, 5).Value = '  -0.17%  '
$ws.Cells.Item(36, 4).Value = '''1.003'
$ws.Cells.Item(36, 5).Value = '  +1.03%  '
$ws.Cells.Item(37, 4).Value = '''0.6245'
$ws.Cells.Item(37, 5).Value = '  +0.64%  '
$ws.Cells.Item(38, 4).Value = '''0.9417'
$ws.Cells.Item(38, 5).Value = '  +3.97%  '
$ws.Cells.Item(39, 4).Value = '''2.045'
$ws.Cells.Item(39, 5).Value = '  +3.47%  '
$ws.Cells.Item(40, 4).Value = '''2.420'
$ws.Cells.Item(40, 5).Value = '  +1.06%  '
$ws.Cells.Item(41, 4).Value = '''0.9996'
$ws.Cells.Item(41, 5).Value = '  -0.06%  '
$ws.Cells.Item(42, 4).Value = '''0.01510'
$ws.Cells.Item(42, 5).Value = '  +1.78%  '
$ws.Cells.Item(43, 4).Value = '''5.575'
$ws.Cells.Item(43, 5).Value = '  +3.34%  '
$ws.Cells.Item(44, 4).Value = '''99.49'
$ws.Cells.Item(44, 5).Value = '  -0.77%  '
$ws.Cells.Item(45, 4).Value = '''0.3861'
$ws.Cells.Item(45, 5).Value = '  +0.56%  '
$ws.Cells.Item(46, 4).Value = '''6.921'
$ws.Cells.Item(46, 5).Value = '  +3.01%  '
$ws.Cells.Item(47, 4).Value = '''0.1160'
$ws.Cells.Item(47, 5).Value = '  +1.17%  '
$ws.Cells.Item(48, 4).Value = '''0.05381'
$ws.Cells.Item(48, 5).Value = '  +0.47%  '
$ws.Cells.Item(49, 4).Value = '''7.888'
$ws.Cells.Item(49, 5).Value = '  +2.76%  '
$ws.Cells.Item(50, 4).Value = '''30.31'
$ws.Cells.Item(50, 5).Value = '  +0.67%  '
$ws.Cells.Item(51, 4).Value = '''51.69'
$ws.Cells.Item(51, 5).Value = '  +1.44%  '
